$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This report is regenerated whenever a file is handed back/re-handed-off.
# In this run, "a9d6742e-7d4c-4504-a071-49a62fa8d0b9" now sorts/lists before
# "1c9b1662-28ba-4c27-9645-463ee3c9dc71" (row 2 vs row 3 swap on every
# sheet), and 1c9b1662's entry has been re-handed-off (new handoff, status
# reset to "Not yet handed off") after its previous handback.
# ---------------------------------------------------------------------------

function Set-LinkDisplay($ws, [string]$addr, [string]$text) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
        }
    }
}

# --- Sheet "Overview" --------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "a9d6742e-7d4c-4504-a071-49a62fa8d0b9.md"
$ws1.Range("B2").Value = "Handed back"
$ws1.Range("C2").Value = "Handed back"

$ws1.Range("A3").Value = "1c9b1662-28ba-4c27-9645-463ee3c9dc71.md"
$ws1.Range("B3").Value = "Not yet handed off"
$ws1.Range("C3").Value = "Not yet handed off"

Set-LinkDisplay $ws1 '$A$2' "a9d6742e-7d4c-4504-a071-49a62fa8d0b9.md"
Set-LinkDisplay $ws1 '$A$3' "1c9b1662-28ba-4c27-9645-463ee3c9dc71.md"

# --- Sheet "zh-cn" -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "a9d6742e-7d4c-4504-a071-49a62fa8d0b9.md"
$ws2.Range("B2").Value = "Handed back"
$ws2.Range("C2").Value = "a9d6742e-7d4c-4504-a071-49a62fa8d0b9.403b532d8323f11c3af5ccf1b83b3ff6487b832a.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-01-11 03:39:51"
$ws2.Range("E2").Value = "a9d6742e-7d4c-4504-a071-49a62fa8d0b9.md"
$ws2.Range("F2").Value = "a9d6742e-7d4c-4504-a071-49a62fa8d0b9.403b532d8323f11c3af5ccf1b83b3ff6487b832a.zh-cn.xlf"
$ws2.Range("G2").Value = "2016-01-11 03:41:02"
$ws2.Range("H2").Value = "Include"

$ws2.Range("A3").Value = "1c9b1662-28ba-4c27-9645-463ee3c9dc71.md"
$ws2.Range("B3").Value = "Not yet handed off"
$ws2.Range("C3").Value = "1c9b1662-28ba-4c27-9645-463ee3c9dc71.3b1d1b2107b976dc0f6861e4da14db54ecf15bd7.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-01-11 03:42:29"
$ws2.Range("E3").Value = "1c9b1662-28ba-4c27-9645-463ee3c9dc71.md"
$ws2.Range("F3").Value = "1c9b1662-28ba-4c27-9645-463ee3c9dc71.3b1d1b2107b976dc0f6861e4da14db54ecf15bd7.zh-cn.xlf"
$ws2.Range("G3").Value = "2016-01-11 03:41:02"
$ws2.Range("H3").Value = "Include"

Set-LinkDisplay $ws2 '$A$2' "a9d6742e-7d4c-4504-a071-49a62fa8d0b9.md"
Set-LinkDisplay $ws2 '$C$2' "a9d6742e-7d4c-4504-a071-49a62fa8d0b9.403b532d8323f11c3af5ccf1b83b3ff6487b832a.zh-cn.xlf"
Set-LinkDisplay $ws2 '$E$2' "a9d6742e-7d4c-4504-a071-49a62fa8d0b9.md"
Set-LinkDisplay $ws2 '$F$2' "a9d6742e-7d4c-4504-a071-49a62fa8d0b9.403b532d8323f11c3af5ccf1b83b3ff6487b832a.zh-cn.xlf"
Set-LinkDisplay $ws2 '$A$3' "1c9b1662-28ba-4c27-9645-463ee3c9dc71.md"
Set-LinkDisplay $ws2 '$C$3' "1c9b1662-28ba-4c27-9645-463ee3c9dc71.3b1d1b2107b976dc0f6861e4da14db54ecf15bd7.zh-cn.xlf"
Set-LinkDisplay $ws2 '$E$3' "1c9b1662-28ba-4c27-9645-463ee3c9dc71.md"
Set-LinkDisplay $ws2 '$F$3' "1c9b1662-28ba-4c27-9645-463ee3c9dc71.3b1d1b2107b976dc0f6861e4da14db54ecf15bd7.zh-cn.xlf"

# --- Sheet "de-de" -------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "a9d6742e-7d4c-4504-a071-49a62fa8d0b9.md"
$ws3.Range("B2").Value = "Handed back"
$ws3.Range("C2").Value = "a9d6742e-7d4c-4504-a071-49a62fa8d0b9.403b532d8323f11c3af5ccf1b83b3ff6487b832a.de-de.xlf"
$ws3.Range("D2").Value = "2016-01-11 03:40:08"
$ws3.Range("E2").Value = "a9d6742e-7d4c-4504-a071-49a62fa8d0b9.md"
$ws3.Range("F2").Value = "a9d6742e-7d4c-4504-a071-49a62fa8d0b9.403b532d8323f11c3af5ccf1b83b3ff6487b832a.de-de.xlf"
$ws3.Range("G2").Value = "2016-01-11 03:41:29"
$ws3.Range("H2").Value = "Include"

$ws3.Range("A3").Value = "1c9b1662-28ba-4c27-9645-463ee3c9dc71.md"
$ws3.Range("B3").Value = "Not yet handed off"
$ws3.Range("C3").Value = "1c9b1662-28ba-4c27-9645-463ee3c9dc71.3b1d1b2107b976dc0f6861e4da14db54ecf15bd7.de-de.xlf"
$ws3.Range("D3").Value = "2016-01-11 03:42:48"
$ws3.Range("E3").Value = "1c9b1662-28ba-4c27-9645-463ee3c9dc71.md"
$ws3.Range("F3").Value = "1c9b1662-28ba-4c27-9645-463ee3c9dc71.3b1d1b2107b976dc0f6861e4da14db54ecf15bd7.de-de.xlf"
$ws3.Range("G3").Value = "2016-01-11 03:41:29"
$ws3.Range("H3").Value = "Include"

Set-LinkDisplay $ws3 '$A$2' "a9d6742e-7d4c-4504-a071-49a62fa8d0b9.md"
Set-LinkDisplay $ws3 '$C$2' "a9d6742e-7d4c-4504-a071-49a62fa8d0b9.403b532d8323f11c3af5ccf1b83b3ff6487b832a.de-de.xlf"
Set-LinkDisplay $ws3 '$E$2' "a9d6742e-7d4c-4504-a071-49a62fa8d0b9.md"
Set-LinkDisplay $ws3 '$F$2' "a9d6742e-7d4c-4504-a071-49a62fa8d0b9.403b532d8323f11c3af5ccf1b83b3ff6487b832a.de-de.xlf"
Set-LinkDisplay $ws3 '$A$3' "1c9b1662-28ba-4c27-9645-463ee3c9dc71.md"
Set-LinkDisplay $ws3 '$C$3' "1c9b1662-28ba-4c27-9645-463ee3c9dc71.3b1d1b2107b976dc0f6861e4da14db54ecf15bd7.de-de.xlf"
Set-LinkDisplay $ws3 '$E$3' "1c9b1662-28ba-4c27-9645-463ee3c9dc71.md"
Set-LinkDisplay $ws3 '$F$3' "1c9b1662-28ba-4c27-9645-463ee3c9dc71.3b1d1b2107b976dc0f6861e4da14db54ecf15bd7.de-de.xlf"
